# Auto-generated update script for cryptos.xlsx price refresh
# Applies the row-level text changes described by the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "0.628") need to be
# forced to Text format first, otherwise Excel auto-converts the assigned
# string into a numeric value, which would change the cell type from text
# to number (the source data keeps these as text/inline strings).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2
$ws.Range('D2').Value = '66.907.62'
$ws.Range('E2').Value = '  +1.87%  '

# Row 3
$ws.Range('D3').Value = '3.438.54'
$ws.Range('E3').Value = '  +1.12%  '

# Row 4
$ws.Range('E4').Value = '  +0.19%  '

# Row 5
Set-TextValue $ws.Range('D5') '575.41'
$ws.Range('E5').Value = '  +2.66%  '

# Row 6
Set-TextValue $ws.Range('D6') '187.23'
$ws.Range('E6').Value = '  +6.28%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.628'
$ws.Range('E7').Value = '  -0.16%  '

# Row 8
$ws.Range('D8').Value = '3.429.97'
$ws.Range('E8').Value = '  +1.21%  '

# Row 9
$ws.Range('E9').Value = '  +0.00%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.171'
$ws.Range('E10').Value = '  -0.98%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.640'
$ws.Range('E11').Value = '  +0.28%  '

# Row 12
Set-TextValue $ws.Range('D12') '57.51'
$ws.Range('E12').Value = '  +6.64%  '

# Row 13
Set-TextValue $ws.Range('D13') '0.0000275'
$ws.Range('E13').Value = '  -1.72%  '

# Row 14
Set-TextValue $ws.Range('D14') '9.42'
$ws.Range('E14').Value = '  +1.96%  '

# Row 15
$ws.Range('D15').Value = '4.003.42'
$ws.Range('E15').Value = '  +1.68%  '

# Row 16
Set-TextValue $ws.Range('D16') '18.87'
$ws.Range('E16').Value = '  +2.88%  '

# Row 17
$ws.Range('D17').Value = '3.446.48'
$ws.Range('E17').Value = '  +1.86%  '

# Row 18
$ws.Range('D18').Value = '66.962.54'
$ws.Range('E18').Value = '  +2.31%  '

# Row 19
$ws.Range('E19').Value = '  -0.67%  '

# Row 20
Set-TextValue $ws.Range('D20') '12.00'
$ws.Range('E20').Value = '  +0.70%  '

# Row 21
Set-TextValue $ws.Range('D21') '1.01'
$ws.Range('E21').Value = '  +1.13%  '

# Row 22
Set-TextValue $ws.Range('D22') '489.11'
$ws.Range('E22').Value = '  +5.20%  '

# Row 23
Set-TextValue $ws.Range('D23') '5.58'
$ws.Range('E23').Value = '  +12.03%  '

# Row 24
Set-TextValue $ws.Range('D24') '16.90'
$ws.Range('E24').Value = '  +17.76%  '

# Row 25
Set-TextValue $ws.Range('D25') '4.34'
$ws.Range('E25').Value = '  +4.82%  '

# Row 26
Set-TextValue $ws.Range('D26') '89.30'
$ws.Range('E26').Value = '  +2.08%  '

# Row 27
Set-TextValue $ws.Range('D27') '2.95'
$ws.Range('E27').Value = '  +0.54%  '

# Row 28
Set-TextValue $ws.Range('D28') '10.91'
$ws.Range('E28').Value = '  +1.45%  '

# Row 29
Set-TextValue $ws.Range('D29') '8.97'
$ws.Range('E29').Value = '  +2.32%  '

# Row 30
Set-TextValue $ws.Range('D30') '31.10'
$ws.Range('E30').Value = '  -0.20%  '

# Row 31
Set-TextValue $ws.Range('D31') '7.35'
$ws.Range('E31').Value = '  +11.75%  '

# Row 32
Set-TextValue $ws.Range('D32') '602.56'
$ws.Range('E32').Value = '  +4.16%  '

# Row 33
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D33') '64.91'
$ws.Range('E33').Value = '  +1.99%  '

# Row 34
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D34') '11.78'
$ws.Range('E34').Value = '  +2.17%  '

# Row 35
Set-TextValue $ws.Range('D35') '0.111'
$ws.Range('E35').Value = '  +2.87%  '

# Row 36
$ws.Range('E36').Value = '  -0.09%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.145'
$ws.Range('E37').Value = '  +1.96%  '

# Row 38
$ws.Range('D38').Value = '0.0₃0775'
$ws.Range('E38').Value = '  +4.14%  '

# Row 39
Set-TextValue $ws.Range('D39') '36.54'
$ws.Range('E39').Value = '  +1.35%  '

# Row 40
Set-TextValue $ws.Range('D40') '0.385'
$ws.Range('E40').Value = '  +2.44%  '

# Row 41
Set-TextValue $ws.Range('D41') '3.44'
$ws.Range('E41').Value = '  -4.44%  '

# Row 42
$ws.Range('D42').Value = '3.183.45'
$ws.Range('E42').Value = '  +2.15%  '

# Row 43
Set-TextValue $ws.Range('D43') '2.87'
$ws.Range('E43').Value = '  +2.24%  '

# Row 44
Set-TextValue $ws.Range('D44') '0.0428'
$ws.Range('E44').Value = '  +2.42%  '

# Row 45
Set-TextValue $ws.Range('D45') '2.55'
$ws.Range('E45').Value = '  +4.15%  '

# Row 46
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range('D46') '3.22'
$ws.Range('E46').Value = '  +2.08%  '

# Row 47
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D47') '0.135'
$ws.Range('E47').Value = '  +0.60%  '

# Row 48
Set-TextValue $ws.Range('D48') '2.64'
$ws.Range('E48').Value = '  +13.64%  '

# Row 49
Set-TextValue $ws.Range('D49') '1.00'
$ws.Range('E49').Value = '  +0.36%  '

# Row 50
Set-TextValue $ws.Range('D50') '8.61'
$ws.Range('E50').Value = '  +1.57%  '

# Row 51
$ws.Range('B51').Value = 'LidoDAOToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D51') '3.17'
$ws.Range('E51').Value = '  +1.15%  '
